$wb = $excel.ActiveWorkbook

# --- Class sheet updates ---
$ws = $wb.Worksheets.Item("Class")
$ws.Activate()

# Update test-data text values (shared strings)
$ws.Range("D3").Value = "Kevin Thomas"
$ws.Range("C2").Value = "Playwrighters team four"
$ws.Range("C3").Value = "updated Playwrighters"

# Resize columns A:C (also clears their "best fit" auto-sizing flag)
$ws.Columns.Item(1).ColumnWidth = 25.166666666666668
$ws.Columns.Item(2).ColumnWidth = 22.333333333333332
$ws.Columns.Item(3).ColumnWidth = 23.166666666666668

# Move the active selection to F5
[void]$ws.Range("F5").Select()
